$d = $word.ActiveDocument
$result = $d.Content.Find.Execute("November 12", $true, $false, $false, $false, $false, $true, 1, $false, "November 19", 2)
Write-Host "Find result: $result"
